# Updates the cryptos list data (Price and Volume(1h) columns, plus two
# row swaps) to match the latest scrape, per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.090.44"
$ws.Range("E2").Value = "  -2.76%  "
$ws.Range("D3").Value = "3.513.36"
$ws.Range("E3").Value = "  -4.55%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "'578.04"
$ws.Range("E5").Value = "  -0.43%  "
$ws.Range("D6").Value = "'171.44"
$ws.Range("E6").Value = "  -3.25%  "
$ws.Range("B7").Value = "XRP"
$ws.Range("C7").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D7").Value = "'0.608"
$ws.Range("E7").Value = "  -0.82%  "
$ws.Range("B8").Value = "LidoStakedEther"
$ws.Range("C8").Value = "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth"
$ws.Range("D8").Value = "3.505.26"
$ws.Range("E8").Value = "  -4.48%  "
$ws.Range("E9").Value = "  -0.07%  "
$ws.Range("E10").Value = "  -5.20%  "
$ws.Range("D11").Value = "'6.55"
$ws.Range("E11").Value = "  -1.05%  "
$ws.Range("E12").Value = "  -4.27%  "
$ws.Range("D13").Value = "'46.84"
$ws.Range("E13").Value = "  -3.98%  "
$ws.Range("D14").Value = "'0.0000273"
$ws.Range("E14").Value = "  -4.26%  "
$ws.Range("D15").Value = "4.090.69"
$ws.Range("E15").Value = "  -4.25%  "
$ws.Range("E16").Value = "  -4.76%  "
$ws.Range("D17").Value = "'622.29"
$ws.Range("E17").Value = "  -8.14%  "
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "3.515.49"
$ws.Range("E18").Value = "  -4.75%  "
$ws.Range("B19").Value = "WrappedBTC"
$ws.Range("C19").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D19").Value = "69.086.02"
$ws.Range("E19").Value = "  -2.92%  "
$ws.Range("E20").Value = "  +0.00%  "
$ws.Range("D21").Value = "'17.43"
$ws.Range("E21").Value = "  -2.42%  "
$ws.Range("D22").Value = "'11.14"
$ws.Range("E22").Value = "  -3.39%  "
$ws.Range("D23").Value = "'0.883"
$ws.Range("E23").Value = "  -5.81%  "
$ws.Range("D24").Value = "'15.94"
$ws.Range("E24").Value = "  -7.73%  "
$ws.Range("D25").Value = "'97.43"
$ws.Range("E25").Value = "  -4.25%  "
$ws.Range("D26").Value = "'3.80"
$ws.Range("E26").Value = "  -4.15%  "
$ws.Range("E27").Value = "  +0.03%  "
$ws.Range("D28").Value = "'2.64"
$ws.Range("E28").Value = "  -6.14%  "
$ws.Range("D29").Value = "'9.33"
$ws.Range("E29").Value = "  -8.59%  "
$ws.Range("D30").Value = "'32.59"
$ws.Range("E30").Value = "  -6.67%  "
$ws.Range("D31").Value = "'3.16"
$ws.Range("E31").Value = "  -7.11%  "
$ws.Range("D32").Value = "'8.53"
$ws.Range("E32").Value = "  -6.70%  "
$ws.Range("D33").Value = "'1.32"
$ws.Range("E33").Value = "  -7.22%  "
$ws.Range("D34").Value = "'7.01"
$ws.Range("E34").Value = "  -5.22%  "
$ws.Range("D35").Value = "'633.98"
$ws.Range("E35").Value = "  +8.51%  "
$ws.Range("D36").Value = "'10.73"
$ws.Range("E36").Value = "  -3.75%  "
$ws.Range("D37").Value = "'0.103"
$ws.Range("E37").Value = "  -5.06%  "
$ws.Range("D38").Value = "'3.41"
$ws.Range("E38").Value = "  -16.15%  "
$ws.Range("D39").Value = "'56.64"
$ws.Range("E39").Value = "  -3.79%  "
$ws.Range("D40").Value = "'1.00"
$ws.Range("E40").Value = "  +0.08%  "
$ws.Range("E41").Value = "  -1.61%  "
$ws.Range("E42").Value = "  -5.33%  "
$ws.Range("D43").Value = "3.370.68"
$ws.Range("E43").Value = "  -8.22%  "
$ws.Range("D44").Value = "'0.327"
$ws.Range("E44").Value = "  -5.88%  "
$ws.Range("D45").Value = "'32.88"
$ws.Range("E45").Value = "  -6.72%  "
$ws.Range("D46").Value = "0.0₃0689"
$ws.Range("E46").Value = "  -9.46%  "
$ws.Range("E47").Value = "  -6.55%  "
$ws.Range("E48").Value = "  -3.63%  "
$ws.Range("D49").Value = "'0.130"
$ws.Range("E49").Value = "  -2.09%  "
$ws.Range("D50").Value = "'132.74"
$ws.Range("E50").Value = "  -2.37%  "
$ws.Range("E51").Value = "  +14.92%  "
